$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ODI Batting Extra")

# Insert 4 new rows at row 2, shifting existing data (rows 2-21) down to rows 6-25
$ws.Rows("2:5").Insert()

# Populate the 4 newly inserted rows with scraped data
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "4248"
$ws.Range("B2").NumberFormat = "General"
$ws.Range("B2").Value = 3
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "1"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "9.78%"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "NO"

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "4249"
$ws.Range("B3").NumberFormat = "General"
$ws.Range("B3").Value = 3
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "1"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "0"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.78%"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "NO"

$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "4437"
$ws.Range("B4").NumberFormat = "General"
$ws.Range("B4").Value = 2
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "3"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "10.93%"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "NO"

$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "4621"
$ws.Range("B5").NumberFormat = "General"
$ws.Range("B5").Value = 2
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "6"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "2"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "20.78%"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "NO"
